$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E width change (4 -> 5)
# NOTE: Excel's ColumnWidth COM property is in "characters" of the default
# font and gets re-quantized to a pixel grid before being stored as the
# raw OOXML <col width> (roughly raw = ColumnWidth + 0.8333...). Feeding it
# 4.166666666666667 round-trips to exactly width="5" in the saved XML.
$ws.Columns.Item(5).ColumnWidth = 4.166666666666667

# Row 1
$ws.Range("E1").Value = 200
$ws.Range("J1").Value = 11.469571800000001
$ws.Range("K1").Value = 0.000056270252276391375
$ws.Range("L1").Value = 0.00000035884963364348802

# Row 2
$ws.Range("E2").Value = 200
$ws.Range("J2").Value = 85.998538100000005
$ws.Range("K2").Value = 0.000057291579179885588
$ws.Range("L2").Value = 0.00000029056295748841403

# Row 3
$ws.Range("E3").Value = 200
$ws.Range("I3").Value = 249
$ws.Range("J3").Value = 87.064598700000005
$ws.Range("K3").Value = 0.00028043565219437738
$ws.Range("L3").Value = -0.000012757075532769695
$ws.Range("N3").Value = 82

# Row 4
$ws.Range("E4").Value = 200
$ws.Range("I4").Value = 269
$ws.Range("J4").Value = 213.07385719999999
$ws.Range("K4").Value = 0.00056235209555466348
$ws.Range("L4").Value = -0.000044311236558517197
$ws.Range("N4").Value = 172

# Row 5
$ws.Range("E5").Value = 400
$ws.Range("J5").Value = 11.617077200000001
$ws.Range("K5").Value = 0.000044798810441504244
$ws.Range("L5").Value = 0.00000022382705932867137

# Row 6
$ws.Range("E6").Value = 400
$ws.Range("J6").Value = 123.7624783
$ws.Range("K6").Value = 0.000052544187840064538
$ws.Range("L6").Value = 0.00000012266450801404607

# Row 7
$ws.Range("E7").Value = 400
$ws.Range("J7").Value = 97.262322100000006
$ws.Range("K7").Value = 0.00030585836905561159
$ws.Range("L7").Value = -0.000027164391021007578

# Row 8
$ws.Range("E8").Value = 400
$ws.Range("J8").Value = 194.27795739999999
$ws.Range("K8").Value = 0.00046422798192669035
$ws.Range("L8").Value = -0.00003819826392043179
$ws.Range("N8").Value = 172

# Row 9
$ws.Range("E9").Value = 600
$ws.Range("J9").Value = 10.9386533
$ws.Range("K9").Value = 0.00030741117852550737
$ws.Range("L9").Value = 0.000001590219940428454

# Row 10
$ws.Range("E10").Value = 600
$ws.Range("I10").Value = 270
$ws.Range("J10").Value = 92.545244400000001
$ws.Range("K10").Value = 0.000052308786186339162
$ws.Range("L10").Value = -0.0000018850709496624299

# Row 11
$ws.Range("E11").Value = 600
$ws.Range("J11").Value = 91.648784699999993
$ws.Range("K11").Value = 0.00025894574966001116
$ws.Range("L11").Value = -0.000019674145135171105

# Row 12
$ws.Range("E12").Value = 600
$ws.Range("I12").Value = 493
$ws.Range("J12").Value = 465.7934128
$ws.Range("K12").Value = 0.0016784108346536986
$ws.Range("L12").Value = -0.0024044526872350658
$ws.Range("M12").Value = 4
$ws.Range("N12").Value = 336

# Row 13
$ws.Range("E13").Value = 800
$ws.Range("J13").Value = 11.8562732
$ws.Range("K13").Value = 0.00025028350163358581
$ws.Range("L13").Value = 0.0000012326467092179433

# Row 14
$ws.Range("E14").Value = 800
$ws.Range("I14").Value = 270
$ws.Range("J14").Value = 113.974687
$ws.Range("K14").Value = 0.000052793127088612835
$ws.Range("L14").Value = -0.0000014276651805418591
$ws.Range("M14").Value = 1
$ws.Range("N14").Value = 53

# Row 15
$ws.Range("E15").Value = 800
$ws.Range("I15").Value = 250
$ws.Range("J15").Value = 92.5555181
$ws.Range("K15").Value = 0.00091414444702575537
$ws.Range("L15").Value = -0.0000091632819965110627
$ws.Range("N15").Value = 82

# Row 16
$ws.Range("E16").Value = 800
$ws.Range("I16").Value = 275
$ws.Range("J16").Value = 232.82780260000001
$ws.Range("K16").Value = 0.00050787581890743283
$ws.Range("L16").Value = -0.00002271357437524097

# Row 17
$ws.Range("E17").Value = 1000
$ws.Range("J17").Value = 13.0947622
$ws.Range("K17").Value = 0.00020984865798445718
$ws.Range("L17").Value = 0.00000068207856344959414

# Row 18
$ws.Range("E18").Value = 1000
$ws.Range("I18").Value = 270
$ws.Range("J18").Value = 130.86597219999999
$ws.Range("K18").Value = 0.000053628834649543933
$ws.Range("L18").Value = -0.000000050041677592651383
$ws.Range("M18").Value = 1
$ws.Range("N18").Value = 53

# Row 19
$ws.Range("E19").Value = 1000
$ws.Range("I19").Value = 254
$ws.Range("J19").Value = 101.75489810000001
$ws.Range("K19").Value = 0.0079417695032171221
$ws.Range("L19").Value = 0.000084237720868004531
$ws.Range("N19").Value = 83

# Row 20
$ws.Range("E20").Value = 1000
$ws.Range("I20").Value = 497
$ws.Range("J20").Value = 341.550275
$ws.Range("K20").Value = 0.00043028847265680348
$ws.Range("L20").Value = -0.012764643082590781
$ws.Range("M20").Value = 4
$ws.Range("N20").Value = 342
